$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @{
    2 = 131.8749747856535
    3 = 13.01410214894575
    4 = 12.50442725944027
    5 = 16.6060245898227
    6 = 22.67297064101339
    7 = 7.170233386188473
    8 = 6.616699233194383
    9 = 21.66345138657745
    10 = 36.01299248933135
    11 = 10.15190710921538
    12 = 2.228504301739868
    13 = 6.029136341286721
    14 = 1.444330918331574
    15 = 2.527957204179294
    16 = 18.0420373719754
    17 = 18.93434652570905
    18 = 18.34073408021675
    19 = 6.483609054332415
    20 = 26.80632812771516
    21 = 69.1902567358651
    22 = 12.73128551886407
    23 = 2.372937393573026
    24 = 22.98603503901824
    25 = 6.618967815788621
    26 = 12.61256302976561
    27 = 28.30434883411037
    28 = 4.909968928129271
    29 = 12.50140248264795
    30 = 2.366887839988391
    31 = 2.440238677202089
    32 = 4.309550734854261
    33 = 5.083137399489445
    34 = 96.80646884972326
    35 = 8.534407719523633
    36 = 22.53534329696295
    37 = 4.202927352925072
    38 = 9.847917041587479
    39 = 8.849740700122727
    40 = 7.673858722109325
    41 = 5.766736954553185
    42 = 270.230251978327
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}

